$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCase1")

$ws.Range("A3").Value = "abcd"
$ws.Range("B3").Value = "hbjbhj"

$ws.Range("B3").Select()
